$d = $word.ActiveDocument

$d.Content.Find.Execute("implementacion", $true, $false, $false, $false, $false,
                         $true, 1, $false, "implementación", 2)

$d.Content.Find.Execute("catalago", $true, $false, $false, $false, $false,
                         $true, 1, $false, "catálogo", 2)
